$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = "Liam (3); Neil (8)"
$ws.Range("I3").Value = "Paul (3); Tara (7)"
$ws.Range("I4").Value = "Anqi (32); Riyansh (10)"
$ws.Range("I5").Value = "Varun (13); Amyra (6)"
$ws.Range("I6").Value = "Hugh (7); Myra (9)"
$ws.Range("I7").Value = "Hugh (7)"
